$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Blad1" -> "Objecten"
# Strip the old Classificatiecode/Waarde/Eigenschappen/"Aanwezig in Project"
# columns (B, C, G, H) and replace the remaining single column with the new
# "Onderdeel" list.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Objecten"

# Drop columns B:H (old Classificatiecode/Waarde/Eigenschappen/Aanwezig in
# Project data + their column widths), shifting everything left so only the
# original column A remains.
$ws1.Range("B1:H9").Delete(-4159) | Out-Null

# Clear out the old column A content (the "Onderdeel" / numbered rows) so we
# can write the new list fresh.
$ws1.Cells.Clear() | Out-Null

$ws1.Columns.Item(1).ColumnWidth = 23.666666666666668

$objecten = @(
    "Onderdeel",
    "Brandmeldinstallatie",
    "Sprinklerinstallatie",
    "Rolluik entree",
    "Rolluik fireshield",
    "Schuifdeur entree",
    "Inbraakbeveiliging",
    "Video bewakingssyteem",
    "Telefooninstallatie"
)
for ($i = 0; $i -lt $objecten.Length; $i++) {
    $row = $i + 1
    $ws1.Range("A$row").Value = $objecten[$i]
}

$ws1.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: new "Parameters" sheet, placed right after "Objecten".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Parameters"

$ws2.Columns.Item(1).ColumnWidth = 22.322916666666668
$ws2.Columns.Item(2).ColumnWidth = 11.166666666666666

$ws2.Range("A1").Value = "Eigenschappen"
$ws2.Range("A2").Value = "Aanwezig in Project"
$ws2.Range("B2").Value = "Boolean"
$ws2.Range("A3").Value = "Standalone"
$ws2.Range("B3").Value = "Dropdown"
$ws2.Range("C3").Value = "Standalone, WC, test"
$ws2.Range("B1").Value = "Type"
$ws2.Range("C1").Value = "Options"

$ws2.Range("A2").Select() | Out-Null
$ws2.Activate() | Out-Null
